$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/7"
$display = "Priya.ranpura by PriyaRanpura · Pull Request #7 · dhavalkeerthi/MRIInterns2026A"

$ws.Hyperlinks.Add($ws.Range("B16"), $target, $null, $null, $target)
$ws.Range("B16").Value = $display

$ws.Range("B16").Select() | Out-Null
